$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert a new sheet right after Sheet1, named "Sheet2"
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# Populate Sheet2 with its content
$ws2.Range("A1").Value = "Hola!"
$ws2.Range("B2").Value = "Azzrael Code "
$ws2.Range("C3").Value = "YouTube "
$ws2.Range("D4").Value = "subs"

# Set Sheet2's selection/active cell
$ws2.Range("B38").Select() | Out-Null

# Switch back to Sheet1, update its selection, and leave it as the active sheet
$ws1.Activate()
$ws1.Range("A18").Select() | Out-Null
